$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write literal text (preserving any leading apostrophe, and
# preventing "2024-01-19"-style strings from being auto-converted to a
# date serial) by temporarily forcing Text format, then stripping the
# format again so the cell keeps the workbook's default (unstyled) look.
# A leading "'" has to be doubled before assignment: Excel's Value setter
# always treats a single leading "'" as a quote-prefix marker (stripped
# from the stored text) regardless of number format, but a second "'"
# survives as a literal character.
function Set-LiteralText {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    if ($text.StartsWith("'")) {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
    $cell.ClearFormats()
}

# ---- Row 54 (No. 51) ----
$ws.Range("A54").Value = 51
$ws.Range("B54").Value = 626402
$ws.Range("C54").Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
Set-LiteralText $ws.Range("D54") "2024-01-19"
Set-LiteralText $ws.Range("E54") "'00026T"
Set-LiteralText $ws.Range("F54") "'241751303000213"
Set-LiteralText $ws.Range("G54") "Pembayaran Belanja Barang Sesuai Surat Tugas Nomor:B.348,B.349/BPPSDM.1/KP.440/I/2024 Tgl.16-1-2024"
Set-LiteralText $ws.Range("H54") "'626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.955.101.AA.001066"
$ws.Range("I54").Value = "IDR"
$ws.Range("J54").Value = 1
$ws.Range("K54").Value = 39160000
$ws.Range("L54").Value = 39160000

# ---- Row 55 (No. 52) ----
$ws.Range("A55").Value = 52
$ws.Range("B55").Value = 626402
$ws.Range("C55").Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
Set-LiteralText $ws.Range("D55") "2024-01-19"
Set-LiteralText $ws.Range("E55") "'00027T"
Set-LiteralText $ws.Range("F55") "'241751303000212"
Set-LiteralText $ws.Range("G55") "Pembayaran Belanja Barang Sesuai Surat Tugas Nomor:B.348,B.349/BPPSDM.1/KP.440/I/2024 Tgl.16-1-2024"
Set-LiteralText $ws.Range("H55") "'626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.956.101.AA.000256"
$ws.Range("I55").Value = "IDR"
$ws.Range("J55").Value = 1
$ws.Range("K55").Value = 33460000
$ws.Range("L55").Value = 33460000

# ---- Row 56 (No. 53) ----
$ws.Range("A56").Value = 53
$ws.Range("B56").Value = 626402
$ws.Range("C56").Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
Set-LiteralText $ws.Range("D56") "2024-01-18"
Set-LiteralText $ws.Range("E56") "'00029T"
Set-LiteralText $ws.Range("F56") "'241751303000204"
Set-LiteralText $ws.Range("G56") "Penggantian Uang Persediaan KKP Untuk keperluan Belanja Barang (BPP 001 SET BRSDMKP)"
Set-LiteralText $ws.Range("H56") "'626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.955.101.AA.001066"
$ws.Range("I56").Value = "IDR"
$ws.Range("J56").Value = 1
$ws.Range("K56").Value = 5618285
$ws.Range("L56").Value = 5618285
